$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: shift the rows 7..14 down to 8..15 (values + formats), bottom-up so
# we never clobber a row before it has been copied away.
# ---------------------------------------------------------------------------
for ($r = 14; $r -ge 7; $r--) {
    $src = $r
    $dst = $r + 1
    $ws.Range("B$src`:F$src").Copy() | Out-Null
    $ws.Range("B$dst`:F$dst").PasteSpecial(-4163) | Out-Null   # xlPasteValues
    $ws.Range("B$src`:F$src").Copy() | Out-Null
    $ws.Range("B$dst`:F$dst").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 2: give row 8 (the row that used to be row 7, "Creazione JSON utente")
# its new yellow highlight plus a start date.
# ---------------------------------------------------------------------------
$ws.Range("B8:F8").Interior.Color = 65535   # yellow FFFFFF00
$ws.Range("B8").NumberFormat = "mm-dd-yy"
$ws.Range("B8").Value = [DateTime]"2017-01-01"

# ---------------------------------------------------------------------------
# Step 3: build the brand-new row 7 (in between "Gestione Routing" and
# "Creazione JSON utente").
# ---------------------------------------------------------------------------
$ws.Range("B7:F7").Interior.Color = 5296274     # green FF92D050
$ws.Range("B7:F7").NumberFormat = "mm-dd-yy"
$ws.Range("D7:F7").WrapText = $true

$ws.Range("B7").Value = [DateTime]"2016-12-30"
$ws.Range("C7").Value = [DateTime]"2017-01-01"
$ws.Range("D7").Value = "gestione v-if sui template e aggiornamento vue"
$ws.Range("E7").Value = "Aggiornate tutte le librerie di vue. V-if sui template funziona solo all'interno di un template container"

$ws.Rows("7:7").RowHeight = 28

# ---------------------------------------------------------------------------
# Step 4: selection matches the authored workbook.
# ---------------------------------------------------------------------------
$ws.Range("B7").Select() | Out-Null
